# csr implementation and testing complete
#
# Adds the new "mstatus" CSR row (row 6) to the spec sheet, matching the
# existing table layout/styling, and updates the sheet's view state
# (zoom + active selection) the same way the author's Excel session left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row of CSR data -------------------------------------------------
$ws.Range("A6").Value = "mstatus"

# "300" is the CSR address and must be stored as text (like the other
# ADDRESS column entries, e.g. "B00"), not auto-converted to a number.
# Toggling the number format to Text around the assignment forces that
# without leaving a stray "number stored as text" quote-prefix marker,
# and we flip it back to General afterwards so no visible formatting
# change is introduced.
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "300"
$ws.Range("B6").NumberFormat = "General"

$ws.Range("C6").Value = "M"
$ws.Range("D6").Value = "R/W"
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = "standard"
$ws.Range("G6").Value = "machine status register (currently only used for testing)"

# --- View state: selection + zoom ----------------------------------------
$null = $ws.Range("D6").Select()
$excel.ActiveWindow.Zoom = 200
